# New feedback system and Reservation changes and fronted small corrections
#
# Appends 6 new reservation rows (17-22) to Sheet1, mirroring the existing
# A:D layout (name / email / date / time).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LiteralText($cell, [string]$text) {
    # Plain `.Value = "2024-04-23"` gets auto-parsed as a date serial by
    # this host (same smart entry Excel itself does), which would also
    # stamp a number-format style onto the cell. The source rows store
    # these look-like-a-date strings as plain text, so round-trip the
    # text through a formula and flatten it back to a static value -
    # that keeps the stored cell a plain string with no style attached.
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 17
$ws.Range("A17").Value = "priyank"
$ws.Range("B17").Value = "pitliyapriyank@gmail.com"
Set-LiteralText $ws.Range("C17") "2024-04-23"
$ws.Range("D17").Value = "20:39"

# Row 18
$ws.Range("A18").Value = "priyank"
$ws.Range("B18").Value = "pitliyapriyank@gmail.com"
Set-LiteralText $ws.Range("C18") "2024-04-24"
$ws.Range("D18").Value = "19:02"

# Row 19
$ws.Range("A19").Value = "gigi"
$ws.Range("B19").Value = "pitliyapriyank@gmail.com"
Set-LiteralText $ws.Range("C19") "2024-04-02"
$ws.Range("D19").Value = "21:08"

# Row 20
$ws.Range("A20").Value = "gigi"
$ws.Range("B20").Value = "pitliyapriyank22@gmail.com"
Set-LiteralText $ws.Range("C20") "2024-04-16"
$ws.Range("D20").Value = "23:20"

# Row 21
$ws.Range("A21").Value = "priyank"
$ws.Range("B21").Value = "pitliyapriyank22@gmail.com"
Set-LiteralText $ws.Range("C21") "2024-04-25"
$ws.Range("D21").Value = "22:33"

# Row 22
$ws.Range("A22").Value = "gigi"
$ws.Range("B22").Value = "pitliyapriyank22@gmail.com"
Set-LiteralText $ws.Range("C22") "2024-04-14"
$ws.Range("D22").Value = "19:34"
